$d = $word.ActiveDocument

$d.Content.Find.Execute("565÷2=282, 1", $true, $true, $false, $false, $false, $true, 1, $false, "587÷8=73, 3", 2) | Out-Null
$d.Content.Find.Execute("682÷4=170, 2", $true, $true, $false, $false, $false, $true, 1, $false, "835÷6=139, 1", 2) | Out-Null
$d.Content.Find.Execute("311÷3=103, 2", $true, $true, $false, $false, $false, $true, 1, $false, "498÷8=62, 2", 2) | Out-Null
$d.Content.Find.Execute("130÷4=32, 2", $true, $true, $false, $false, $false, $true, 1, $false, "150÷3=50, 0", 2) | Out-Null
$d.Content.Find.Execute("521÷2=260, 1", $true, $true, $false, $false, $false, $true, 1, $false, "589÷5=117, 4", 2) | Out-Null
$d.Content.Find.Execute("124÷5=24, 4", $true, $true, $false, $false, $false, $true, 1, $false, "400÷9=44, 4", 2) | Out-Null
$d.Content.Find.Execute("161÷4=40, 1", $true, $true, $false, $false, $false, $true, 1, $false, "123÷8=15, 3", 2) | Out-Null
$d.Content.Find.Execute("791÷5=158, 1", $true, $true, $false, $false, $false, $true, 1, $false, "235÷5=47, 0", 2) | Out-Null
$d.Content.Find.Execute("385÷6=64, 1", $true, $true, $false, $false, $false, $true, 1, $false, "310÷2=155, 0", 2) | Out-Null
$d.Content.Find.Execute("951÷7=135, 6", $true, $true, $false, $false, $false, $true, 1, $false, "583÷7=83, 2", 2) | Out-Null
$d.Content.Find.Execute("652÷3=217, 1", $true, $true, $false, $false, $false, $true, 1, $false, "516÷5=103, 1", 2) | Out-Null
$d.Content.Find.Execute("362÷4=90, 2", $true, $true, $false, $false, $false, $true, 1, $false, "151÷6=25, 1", 2) | Out-Null
$d.Content.Find.Execute("541÷4=135, 1", $true, $true, $false, $false, $false, $true, 1, $false, "482÷3=160, 2", 2) | Out-Null
$d.Content.Find.Execute("491÷5=98, 1", $true, $true, $false, $false, $false, $true, 1, $false, "102÷5=20, 2", 2) | Out-Null
$d.Content.Find.Execute("572÷6=95, 2", $true, $true, $false, $false, $false, $true, 1, $false, "420÷2=210, 0", 2) | Out-Null
$d.Content.Find.Execute("725÷5=145, 0", $true, $true, $false, $false, $false, $true, 1, $false, "229÷6=38, 1", 2) | Out-Null
$d.Content.Find.Execute("293÷6=48, 5", $true, $true, $false, $false, $false, $true, 1, $false, "415÷4=103, 3", 2) | Out-Null
$d.Content.Find.Execute("776÷5=155, 1", $true, $true, $false, $false, $false, $true, 1, $false, "191÷9=21, 2", 2) | Out-Null
$d.Content.Find.Execute("854÷5=170, 4", $true, $true, $false, $false, $false, $true, 1, $false, "207÷3=69, 0", 2) | Out-Null
$d.Content.Find.Execute("267÷5=53, 2", $true, $true, $false, $false, $false, $true, 1, $false, "726÷5=145, 1", 2) | Out-Null
$d.Content.Find.Execute("100÷6=16, 4", $true, $true, $false, $false, $false, $true, 1, $false, "839÷9=93, 2", 2) | Out-Null
$d.Content.Find.Execute("647÷5=129, 2", $true, $true, $false, $false, $false, $true, 1, $false, "345÷7=49, 2", 2) | Out-Null
$d.Content.Find.Execute("767÷6=127, 5", $true, $true, $false, $false, $false, $true, 1, $false, "562÷9=62, 4", 2) | Out-Null
$d.Content.Find.Execute("691÷5=138, 1", $true, $true, $false, $false, $false, $true, 1, $false, "779÷2=389, 1", 2) | Out-Null
$d.Content.Find.Execute("240÷9=26, 6", $true, $true, $false, $false, $false, $true, 1, $false, "631÷9=70, 1", 2) | Out-Null
